$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "INTEND.DATE"
$ws.Range("D1").Value = "TAX.INTEREST.TYPE:1"
$ws.Range("E1").Value = "CHARGE.CODE:1"
$ws.Range("F1").Value = "CHARGE.AMOUNT:1"
$ws.Range("G1").Value = "DRAWDOWN.ACCOUNT"
$ws.Range("H1").Value = "PRIN.LIQ.ACCT"
$ws.Range("I1").Value = "INT.LIQ.ACCT"
$ws.Range("J1").Value = "CHRG.LIQ.ACCT"
$ws.Range("K1").Value = "FINAL.MATURITY"
$ws.Range("L1").Value = "EXP.DATE"

$ws.Range("L1").Select()
